$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 17 (CreateDate): type changes from DATE to TIMESTAMP, length cleared
$ws.Range("D17").Value = "TIMESTAMP"
$ws.Range("E17").Value = $null

# Row 19 (LastUpdate): type changes from DATE to TIMESTAMP, length cleared
$ws.Range("D19").Value = "TIMESTAMP"
$ws.Range("E19").Value = $null

# Leave the cursor where the author left it when the workbook was saved
$ws.Range("C22").Select()
